$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the stray trailing newline from the "Jack in the Box" JSON example (row 5, column F)
$ws.Range("F5").Value = '{"RNAME":"Jack in the Box","ADDRESS":"San Fransisco","PHONE":"500-004-3003"}'

# Fix up the JSON example text for the "/menu" (POST) row (row 8, column F)
$ws.Range("F8").Value = '{"MNAME":"Dinner","MDETAILS":"All dishes relating to dinner before 9:00PM are stored here","RID":"1"}'

# Fix up the JSON example text for the "/menuItem" (POST) row (row 11, column F), renaming Prime RIB -> PrimeRib
$ws.Range("F11").Value = '{"MITEMNAME":"PrimeRib","MITEMDETAILS":"Burger","MITEMPRICE":4.95,"MID":"1","RID":"1"}'

# Shorter text now wraps onto fewer lines, so these rows shrink from 4 lines to 3
$ws.Rows(5).RowHeight = 43.2
$ws.Rows(8).RowHeight = 43.2
$ws.Rows(11).RowHeight = 43.2

# Move the selection to a single cell (F9) instead of the whole used range
$ws.Range("F9").Select()
